$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top, pushing existing data down.
$ws.Rows.Item(1).Insert()

# Fill in the new header row values.
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"
$ws.Range("C1").Value = "LMRole"
$ws.Range("D1").Value = "TMRole"
$ws.Range("F1").Value = "TimeZone"
$ws.Range("E1").Value = "User language"

# Make the header row bold.
$ws.Range("A1:F1").Font.Bold = $true
